$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks up front; they will be re-created after the
# cell values below have been shifted into their new rows.
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value2 = '2025-10-23 12:44:20'
$ws.Range("B2").Value2 = '【26年5月/17日間/対面】Python Webアプリ開発 研修講師募集(カリキュラム設計含む)'
$ws.Range("C2").Value2 = 'システム開発'
$ws.Range("D2").Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E2").Value2 = '期限情報なし'
$ws.Range("F2").Value2 = 'https://www.lancers.jp/work/detail/5419191'
$ws.Range("G2").Value2 = 295
$ws.Range("H2").Value2 = '🔥Python ◆開発 ◇アプリ'

$ws.Range("A3").Value2 = '2025-10-23 12:44:20'
$ws.Range("B3").Value2 = '【急募】APIを利用した診断サイト構築のフリーランス募集'
$ws.Range("C3").Value2 = 'システム開発'
$ws.Range("D3").Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E3").Value2 = '期限情報なし'
$ws.Range("F3").Value2 = 'https://www.lancers.jp/work/detail/5418643'
$ws.Range("G3").Value2 = 220
$ws.Range("H3").Value2 = '🔥API ◇サイト'

$ws.Range("A4").Value2 = '2025-10-23 12:44:20'
$ws.Range("B4").Value2 = '【急募】施行主向け建築資材配達アプリ開発者を募集します'
$ws.Range("C4").Value2 = 'システム開発'
$ws.Range("D4").Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value2 = '期限情報なし'
$ws.Range("F4").Value2 = 'https://www.lancers.jp/work/detail/5418447'
$ws.Range("G4").Value2 = 100
$ws.Range("H4").Value2 = '◆開発 ◇アプリ'

$ws.Range("A5").Value2 = '2025-10-23 12:44:20'
$ws.Range("B5").Value2 = '【高品質な恋愛マッチングアプリ制作】エンジニア募集'
$ws.Range("C5").Value2 = 'システム開発'
$ws.Range("D5").Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E5").Value2 = '期限情報なし'
$ws.Range("F5").Value2 = 'https://www.lancers.jp/work/detail/5418455'
$ws.Range("G5").Value2 = 45
$ws.Range("H5").Value2 = '◇アプリ'

$ws.Range("A6").Value2 = '2025-10-23 12:44:20'
$ws.Range("B6").Value2 = '【Webarena suiteX/DNS】ドメイン設定変更によるウェブサイト分割とサイト切り替え'
$ws.Range("C6").Value2 = 'システム開発'
$ws.Range("D6").Value2 = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E6").Value2 = '期限情報なし'
$ws.Range("F6").Value2 = 'https://www.lancers.jp/work/detail/5417544'
$ws.Range("G6").Value2 = 30
$ws.Range("H6").Value2 = '◇サイト'

$ws.Range("A7").Value2 = '2025-10-23 12:44:20'
$ws.Range("B7").Value2 = 'ERPシステムの第三者技術検証・品質評価報告書作成'
$ws.Range("C7").Value2 = 'システム開発'
$ws.Range("D7").Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E7").Value2 = '期限情報なし'
$ws.Range("F7").Value2 = 'https://www.lancers.jp/work/detail/5418891'
$ws.Range("G7").Value2 = 40

$ws.Range("A8").Value2 = '2025-10-23 12:44:20'
$ws.Range("B8").Value2 = '【急募】セッション体験を再現するクローンシステム構築依頼'
$ws.Range("C8").Value2 = 'システム開発'
$ws.Range("D8").Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E8").Value2 = '期限情報なし'
$ws.Range("F8").Value2 = 'https://www.lancers.jp/work/detail/5418644'
$ws.Range("G8").Value2 = 40

$ws.Range("A9").Value2 = '2025-10-23 12:44:20'
$ws.Range("B9").Value2 = '【急募】既存の予約システムの料金修正を依頼します'
$ws.Range("C9").Value2 = 'システム開発'
$ws.Range("D9").Value2 = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E9").Value2 = '期限情報なし'
$ws.Range("F9").Value2 = 'https://www.lancers.jp/work/detail/5418759'
$ws.Range("G9").Value2 = 25

$ws.Range("A10").Value2 = '2025-10-23 12:44:20'
$ws.Range("B10").Value2 = 'Stable Diffusion LoRA制作依頼 画風再現+キャラLoRA量産テンプレ構築'
$ws.Range("C10").Value2 = 'システム開発'
$ws.Range("D10").Value2 = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E10").Value2 = '期限情報なし'
$ws.Range("F10").Value2 = 'https://www.lancers.jp/work/detail/5418738'
$ws.Range("G10").Value2 = 18

$ws.Range("A11").Value2 = '2025-10-23 12:44:20'
$ws.Range("B11").Value2 = '【メールマーケティング】戦略立案・実行者募集'
$ws.Range("C11").Value2 = 'システム開発'
$ws.Range("D11").Value2 = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E11").Value2 = '期限情報なし'
$ws.Range("F11").Value2 = 'https://www.lancers.jp/work/detail/5418443'
$ws.Range("G11").Value2 = 18

$ws.Range("A12").Value2 = '2025-10-23 12:44:20'
$ws.Range("B12").Value2 = '【急募】HPの微修正をお手伝いしてくれる方募集!'
$ws.Range("C12").Value2 = 'システム開発'
$ws.Range("D12").Value2 = '5,000 円 ~'
$ws.Range("E12").Value2 = '期限情報なし'
$ws.Range("F12").Value2 = 'https://www.lancers.jp/work/detail/5418445'
$ws.Range("G12").Value2 = 10

# Re-create the URL hyperlinks for column F (rows 2-12), restoring the
# "Hyperlink" cell style that Excel applies to linked cells.
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5419191')
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5418643')
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5418447')
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5418455')
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5417544')
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5418891')
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5418644')
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5418759')
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5418738')
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5418443')
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5418445')
$ws.Range("F12").Style = "Hyperlink"

# Column widths: B grows from 50 to 51 characters, H grows from 12 to 18.
$ws.Columns.Item(2).ColumnWidth = 50.17
$ws.Columns.Item(8).ColumnWidth = 17.17
